$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: replace "HC-12" / 3.8 with "LoRa Module" / "?"
$ws.Range("A5").Value = "LoRa Module"
$ws.Range("B5").Value = "?"

# Update selection to F10
$ws.Range("F10").Select()
